$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (case fix: "Cluster Name" -> "Cluster name", "Active Cases" -> "Active cases")
$ws.Range("A1").Value = "Cluster name"
$ws.Range("B1").Value = "Active cases"

# Update data rows (cluster names + active case counts)
$ws.Range("A2").Value = "3398 BlueCross Elly Kay Mordialloc"
$ws.Range("B2").Value = 32
$ws.Range("A3").Value = "3564 Waverley Valley Aged Care Glen Waverley"
$ws.Range("B3").Value = 12
$ws.Range("A4").Value = "3601 Baptcare Westhaven community"
$ws.Range("B4").Value = 13
$ws.Range("A5").Value = "3646 Mornington Bay Care Community Mount Martha"
$ws.Range("B5").Value = 17
$ws.Range("A6").Value = "3647 Aurrum Aged Care Reservoir"
$ws.Range("B6").Value = 12
$ws.Range("A7").Value = "3653 Fronditha Thalpori St Albans Aged Care"
$ws.Range("B7").Value = 20
$ws.Range("A8").Value = "3975 Aurrum Aged Care Brunswick West"
$ws.Range("B8").Value = 11
$ws.Range("A9").Value = "4257 BlueCross The Gables Camberwell"
$ws.Range("B9").Value = 18
$ws.Range("A10").Value = "4295 Hope Aged Care Sunshine West"
$ws.Range("B10").Value = 16
$ws.Range("A11").Value = "4314 Estia Health Ardeer"
$ws.Range("B11").Value = 19
$ws.Range("A12").Value = "44304 Brighton Primary School Brighton"
$ws.Range("B12").Value = 17
$ws.Range("A13").Value = "44380 Plenty Parklands Primary School Mill Park"
$ws.Range("B13").Value = 10
$ws.Range("A14").Value = "44414 Buangor Primary School Buangor"
$ws.Range("B14").Value = 21
$ws.Range("A15").Value = "44490 Armadale Primary School Armadale"
$ws.Range("B15").Value = 17
$ws.Range("A16").Value = "44584 Badger Creek Primary School Badger Creek"
$ws.Range("B16").Value = 13
$ws.Range("A17").Value = "44593 Torquay P-6 College Torquay"
$ws.Range("B17").Value = 19
$ws.Range("A18").Value = "44761 Coburg North Primary School Coburg"
$ws.Range("B18").Value = 21
$ws.Range("A19").Value = "44828 Cheltenham East Primary School Cheltenham"
$ws.Range("B19").Value = 11
$ws.Range("A20").Value = "44930 St Albans Heights Primary School"
$ws.Range("B20").Value = 13
$ws.Range("A21").Value = "44979 Campbellfield Heights Primary School Campbellfield"
$ws.Range("B21").Value = 18
$ws.Range("A22").Value = "45248 Brookside P-9 College Caroline Springs"
$ws.Range("B22").Value = 14
$ws.Range("A23").Value = "45257 Roxburgh Rise Primary School Roxburgh Park"
$ws.Range("B23").Value = 13
$ws.Range("A24").Value = "45708 St Augustine's Primary School Yarraville"
$ws.Range("B24").Value = 10
$ws.Range("A25").Value = "45719 St Joseph's Primary School Numurkah"
$ws.Range("B25").Value = 14
$ws.Range("A26").Value = "4574 Village Glen Aged Care Residences Mornington"
$ws.Range("B26").Value = 10
$ws.Range("A27").Value = "45861 St Oliver Plunkett Primary School Pascoe Vale"
$ws.Range("B27").Value = 12
$ws.Range("A28").Value = "45902 Mother of God Primary School Ardeer"
$ws.Range("B28").Value = 11
$ws.Range("A29").Value = "45975 St Thomas More Primary School Hadfield"
$ws.Range("B29").Value = 12
$ws.Range("A30").Value = "45988 St Macartan's Parish Primary School"
$ws.Range("B30").Value = 11
$ws.Range("A31").Value = "46078 Corpus Christi Primary School Werribee"
$ws.Range("B31").Value = 13
$ws.Range("A32").Value = "46116 Saint Francis Of Assisi Primary School Baranduda"
$ws.Range("B32").Value = 10
$ws.Range("A33").Value = "46277 St Paul's Anglican Grammar School Warragul"
$ws.Range("B33").Value = 11
$ws.Range("A34").Value = "46306 King's College Warrnambool"
$ws.Range("B34").Value = 11
$ws.Range("A35").Value = "50722 Chairo Christian School Leongatha Campus Leongatha"
$ws.Range("B35").Value = 10
$ws.Range("A36").Value = "51529 Sirius College Primary School Dallas"
$ws.Range("B36").Value = 11
$ws.Range("A37").Value = "Berwick Chase Primary School Berwick"
$ws.Range("B37").Value = 10
$ws.Range("A38").Value = "Kororoit Creek Primary School Burnside Heights Oct-Dec"
$ws.Range("B38").Value = 24
$ws.Range("A39").Value = "Mambourin Enterprises Allara Deer Park"
$ws.Range("B39").Value = 28
$ws.Range("A40").Value = "North St Kilda Children's Centre St Kilda"
$ws.Range("B40").Value = 10
$ws.Range("A41").Value = "Oakleigh South Primary School Oakleigh South"
$ws.Range("B41").Value = 13
$ws.Range("A42").Value = "Rosebud Primary School Rosebud"
$ws.Range("B42").Value = 17
$ws.Range("A43").Value = "Springside Primary School Caroline Springs Nov"
$ws.Range("B43").Value = 15
$ws.Range("A44").Value = "St Christophers Primary School Airport West"
$ws.Range("B44").Value = 13
$ws.Range("A45").Value = "St Louis de Montfort's School Aspendale"
$ws.Range("B45").Value = 17
$ws.Range("A46").Value = "Torquay Hotel Torquay"
$ws.Range("B46").Value = 20
